# Update "想去人数" (column F) values across the "展览", "演出" and "全部类型"
# worksheets to reflect refreshed counts from the upstream data source.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 77
$ws.Range("F3").Value  = 167
$ws.Range("F5").Value  = 51
$ws.Range("F6").Value  = 2743
$ws.Range("F8").Value  = 1628
$ws.Range("F11").Value = 7624
$ws.Range("F15").Value = 6131
$ws.Range("F16").Value = 3250
$ws.Range("F17").Value = 3623
$ws.Range("F24").Value = 282
$ws.Range("F25").Value = 282
$ws.Range("F26").Value = 3615
$ws.Range("F28").Value = 338
$ws.Range("F31").Value = 1082
$ws.Range("F34").Value = 2604
$ws.Range("F35").Value = 1452
$ws.Range("F39").Value = 3241
$ws.Range("F40").Value = 153
$ws.Range("F45").Value = 1268
$ws.Range("F48").Value = 586

# ---- Sheet: 演出 (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value  = 19
$ws.Range("F8").Value  = 39
$ws.Range("F9").Value  = 399
$ws.Range("F18").Value = 12

# ---- Sheet: 全部类型 (All types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 77
$ws.Range("F5").Value  = 167
$ws.Range("F7").Value  = 51
$ws.Range("F9").Value  = 2743
$ws.Range("F10").Value = 1628
$ws.Range("F14").Value = 7624
$ws.Range("F17").Value = 6131
$ws.Range("F18").Value = 3250
$ws.Range("F19").Value = 3623
$ws.Range("F26").Value = 282
$ws.Range("F28").Value = 282
$ws.Range("F29").Value = 3615
$ws.Range("F34").Value = 338
$ws.Range("F38").Value = 2604
$ws.Range("F39").Value = 1452
$ws.Range("F43").Value = 3241
$ws.Range("F47").Value = 1268
$ws.Range("F49").Value = 586
